$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "2022-Q1" sheet, positioned right before the "总计" sheet.
#    We duplicate the "2021-Q4" sheet (same column layout/styles) so the new
#    sheet inherits the correct header/body formatting (bold+border style on
#    row 1 and on column A), then we rename it and overwrite every value.
# ---------------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item(3)   # "2021-Q4"
$totalSheet    = $wb.Worksheets.Item(4)   # "总计" (new sheet goes right before it)

$templateSheet.Copy($totalSheet)

$newSheet = $wb.Worksheets.Item(4)
$newSheet.Name = "2022-Q1"

# The template only has 6 rows (1 header + 5 data); we need 10 rows (1 header
# + 9 data). Extend the column-A / row styling down to rows 7:10 by copying
# the format of the last existing data row.
$newSheet.Range("A6:H6").Copy()
$newSheet.Range("A7:H10").PasteSpecial(-4122)   # xlPasteFormats

# Columns B-G hold text values (fund codes, names and numbers-as-text such as
# "23.23"); force a Text number format so they are not auto-coerced into
# numeric cells. Columns A and H stay numeric.
$newSheet.Range("B2:G10").NumberFormat = "@"

# Row 1 - header
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Row 2
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "005821"
$newSheet.Range("C2").Value = "万家新机遇龙头企业灵活配置混合"
$newSheet.Range("D2").Value = "23.23"
$newSheet.Range("E2").Value = "56.20"
$newSheet.Range("F2").Value = "3.54"
$newSheet.Range("G2").Value = "0.8223"
$newSheet.Range("H2").Value = 3

# Row 3
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "013960"
$newSheet.Range("C3").Value = "万家新机遇成长一年持有期混合A"
$newSheet.Range("D3").Value = "13.29"
$newSheet.Range("E3").Value = "49.51"
$newSheet.Range("F3").Value = "3.23"
$newSheet.Range("G3").Value = "0.4293"
$newSheet.Range("H3").Value = 4

# Row 4
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "013961"
$newSheet.Range("C4").Value = "万家新机遇成长一年持有期混合C"
$newSheet.Range("D4").Value = "3.13"
$newSheet.Range("E4").Value = "49.51"
$newSheet.Range("F4").Value = "3.23"
$newSheet.Range("G4").Value = "0.1011"
$newSheet.Range("H4").Value = 4

# Row 5
$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "011071"
$newSheet.Range("C5").Value = "鹏华安悦一年持有期混合A"
$newSheet.Range("D5").Value = "9.16"
$newSheet.Range("E5").Value = "21.81"
$newSheet.Range("F5").Value = "0.73"
$newSheet.Range("G5").Value = "0.0669"
$newSheet.Range("H5").Value = 3

# Row 6
$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "001067"
$newSheet.Range("C6").Value = "鹏华弘盛灵活配置混合A"
$newSheet.Range("D6").Value = "6.22"
$newSheet.Range("E6").Value = "20.46"
$newSheet.Range("F6").Value = "0.72"
$newSheet.Range("G6").Value = "0.0448"
$newSheet.Range("H6").Value = 3

# Row 7
$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "009232"
$newSheet.Range("C7").Value = "鹏华安惠混合A"
$newSheet.Range("D7").Value = "3.74"
$newSheet.Range("E7").Value = "21.93"
$newSheet.Range("F7").Value = "0.78"
$newSheet.Range("G7").Value = "0.0292"
$newSheet.Range("H7").Value = 3

# Row 8
$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").Value = "001380"
$newSheet.Range("C8").Value = "鹏华弘盛灵活配置混合C"
$newSheet.Range("D8").Value = "0.53"
$newSheet.Range("E8").Value = "20.46"
$newSheet.Range("F8").Value = "0.72"
$newSheet.Range("G8").Value = "0.0038"
$newSheet.Range("H8").Value = 3

# Row 9
$newSheet.Range("A9").Value = 7
$newSheet.Range("B9").Value = "009233"
$newSheet.Range("C9").Value = "鹏华安惠混合C"
$newSheet.Range("D9").Value = "0.32"
$newSheet.Range("E9").Value = "21.93"
$newSheet.Range("F9").Value = "0.78"
$newSheet.Range("G9").Value = "0.0025"
$newSheet.Range("H9").Value = 3

# Row 10
$newSheet.Range("A10").Value = 8
$newSheet.Range("B10").Value = "005493"
$newSheet.Range("C10").Value = "鑫元价值精选灵活配置混合A"
$newSheet.Range("D10").Value = "0.07"
$newSheet.Range("E10").Value = "21.81"
$newSheet.Range("F10").Value = "0.73"
$newSheet.Range("G10").Value = "0.0005"
$newSheet.Range("H10").Value = 3

# ---------------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: add a new first data row for 2022-Q1
#    and shift the previously-existing rows down by one, renumbering the
#    index column (A) accordingly.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(5)   # "总计" is now the 5th sheet

# Row 5 is brand new (the sheet previously only went down to row 4); copy the
# row-4 formatting down first so the A-column keeps its bold/bordered style.
$totalSheet.Range("A4:D4").Copy()
$totalSheet.Range("A5:D5").PasteSpecial(-4122)   # xlPasteFormats

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 9
$totalSheet.Range("D2").Value = 1.5

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 5
$totalSheet.Range("D3").Value = 0.32

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q3"
$totalSheet.Range("C4").Value = 1
$totalSheet.Range("D4").Value = 0.03

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2021-Q2"
$totalSheet.Range("C5").Value = 1
$totalSheet.Range("D5").Value = 0.02
